$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 214, shifting existing
# rows 214:223 down to 216:225 (same as rows 214/215 being duplicated
# with a newer date in the source weekly update).
$ws.Range("A214:A215").EntireRow.Insert()

# New row 214 - "Primera" quality entry for the new week (44509)
$ws.Range("A214").Value = 4
$ws.Range("B214").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C214").Value = "Los Lagos"
$ws.Range("D214").Value = 44509
$ws.Range("E214").Value = 10
$ws.Range("F214").Value = 100112023
$ws.Range("G214").Value = "Brócoli"
$ws.Range("H214").Value = "Sin especificar"
$ws.Range("I214").Value = "Primera"
$ws.Range("J214").Value = 700
$ws.Range("K214").Value = 1200
$ws.Range("L214").Value = 1200
$ws.Range("M214").Value = 1200
$ws.Range("N214").Value = "$/unidad"
$ws.Range("O214").Value = "Región Metropolitana"
$ws.Range("P214").Value = 1200
$ws.Range("Q214").Value = 1
$ws.Range("R214").Value = "Hortaliza"

# New row 215 - "Segunda" quality entry for the new week (44509)
$ws.Range("A215").Value = 4
$ws.Range("B215").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C215").Value = "Los Lagos"
$ws.Range("D215").Value = 44509
$ws.Range("E215").Value = 10
$ws.Range("F215").Value = 100112023
$ws.Range("G215").Value = "Brócoli"
$ws.Range("H215").Value = "Sin especificar"
$ws.Range("I215").Value = "Segunda"
$ws.Range("J215").Value = 700
$ws.Range("K215").Value = 1000
$ws.Range("L215").Value = 1000
$ws.Range("M215").Value = 1000
$ws.Range("N215").Value = "$/unidad"
$ws.Range("O215").Value = "Región Metropolitana"
$ws.Range("P215").Value = 1000
$ws.Range("Q215").Value = 1
$ws.Range("R215").Value = "Hortaliza"
